$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")

# Row 11
$ws.Range("H11").Value = 60.058823
$ws.Range("I11").Value = 60.058823
$ws.Range("K11").Value = 60.058823
$ws.Range("M11").Value = 79.94117700000001

# Row 100
$ws.Range("H100").Value = 1399.5
$ws.Range("I100").Value = 979.6
$ws.Range("K100").Value = 979.6
$ws.Range("M100").Value = -438.6

# Row 137
$ws.Range("H137").Value = 2077
$ws.Range("I137").Value = 2111.0833
$ws.Range("K137").Value = 6333.249899999999
$ws.Range("M137").Value = -3783.249899999999

# Row 138
$ws.Range("H138").Value = 7222.9414
$ws.Range("I138").Value = 1066.1666
$ws.Range("J138").Value = 21999.2
$ws.Range("K138").Value = 3198.4998
$ws.Range("L138").Value = 65997.60000000001
$ws.Range("M138").Value = 1941.5002
$ws.Range("N138").Value = -76277.60000000001


$ws = $wb.Worksheets.Item("ARM")

# Row 25
$ws.Range("H25").Value = 329.5
$ws.Range("I25").Value = 172.66667
$ws.Range("J25").Value = 800
$ws.Range("K25").Value = 172.66667
$ws.Range("L25").Value = 800
$ws.Range("M25").Value = 229.33333
$ws.Range("N25").Value = -1604

# Row 29
$ws.Range("H29").Value = 1000
$ws.Range("I29").Value = 1000
$ws.Range("J29").Value = 1000
$ws.Range("K29").Value = 1000
$ws.Range("L29").Value = 1000
$ws.Range("M29").Value = -692
$ws.Range("N29").Value = -1616

# Row 32
$ws.Range("H32").Value = 4880.4316
$ws.Range("I32").Value = 3733.7097
$ws.Range("K32").Value = 3733.7097
$ws.Range("M32").Value = -3446.7097

# Row 41
$ws.Range("H41").Value = 15000
$ws.Range("I41").Value = 15000
$ws.Range("K41").Value = 15000
$ws.Range("M41").Value = -14586


$ws = $wb.Worksheets.Item("BSM")

# Row 12
$ws.Range("H12").Value = 4083.1667
$ws.Range("I12").Value = 4750
$ws.Range("J12").Value = 3749.75
$ws.Range("K12").Value = 4750
$ws.Range("L12").Value = 3749.75
$ws.Range("M12").Value = -4582
$ws.Range("N12").Value = -4085.75

# Row 20
$ws.Range("H20").Value = 6132.5835
$ws.Range("I20").Value = 6776
$ws.Range("K20").Value = 6776
$ws.Range("M20").Value = -6529

# Row 26
$ws.Range("H26").Value = 34333.332
$ws.Range("I26").Value = 34333.332
$ws.Range("K26").Value = 34333.332
$ws.Range("M26").Value = -34041.332

# Row 86
$ws.Range("H86").Value = 3683
$ws.Range("I86").Value = 1839.8
$ws.Range("K86").Value = 1839.8
$ws.Range("M86").Value = -716.8

# Row 89
$ws.Range("H89").Value = 3683
$ws.Range("I89").Value = 1839.8
$ws.Range("K89").Value = 9199
$ws.Range("M89").Value = -3583

# Row 96
$ws.Range("H96").Value = 10679.25
$ws.Range("I96").Value = 10679.25
$ws.Range("K96").Value = 10679.25
$ws.Range("M96").Value = -7933.25


$ws = $wb.Worksheets.Item("CRP")

# Row 22
$ws.Range("H22").Value = 312.5
$ws.Range("I22").Value = 295
$ws.Range("J22").Value = 400
$ws.Range("K22").Value = 295
$ws.Range("L22").Value = 400
$ws.Range("M22").Value = 55
$ws.Range("N22").Value = -1100

# Row 32
$ws.Range("H32").Value = 77315.30499999999
$ws.Range("I32").Value = 83658.25
$ws.Range("J32").Value = 1200
$ws.Range("K32").Value = 83658.25
$ws.Range("L32").Value = 1200
$ws.Range("M32").Value = -83342.25
$ws.Range("N32").Value = -1832

# Row 86
$ws.Range("H86").Value = 9380.833000000001
$ws.Range("I86").Value = 8067.857
$ws.Range("J86").Value = 11219
$ws.Range("K86").Value = 8067.857
$ws.Range("L86").Value = 11219
$ws.Range("M86").Value = -6944.857
$ws.Range("N86").Value = -13465

# Row 89
$ws.Range("H89").Value = 9380.833000000001
$ws.Range("I89").Value = 8067.857
$ws.Range("J89").Value = 11219
$ws.Range("K89").Value = 40339.285
$ws.Range("L89").Value = 56095
$ws.Range("M89").Value = -34723.285
$ws.Range("N89").Value = -67327

# Row 107
$ws.Range("H107").Value = 1037.7693
$ws.Range("I107").Value = 802.13336
$ws.Range("K107").Value = 802.13336
$ws.Range("M107").Value = 1117.86664

# Row 122
$ws.Range("H122").Value = 3063
$ws.Range("I122").Value = 3197.8462
$ws.Range("J122").Value = 2624.75
$ws.Range("K122").Value = 9593.5386
$ws.Range("L122").Value = 7874.25
$ws.Range("M122").Value = -7143.5386
$ws.Range("N122").Value = -12774.25

# Row 132
$ws.Range("H132").Value = 1812.0212
$ws.Range("I132").Value = 1485.3489
$ws.Range("K132").Value = 4456.0467
$ws.Range("M132").Value = -1926.0467

# Row 134
$ws.Range("H134").Value = 1894.3962
$ws.Range("I134").Value = 1616.1621
$ws.Range("K134").Value = 4848.4863
$ws.Range("M134").Value = -2313.4863


$ws = $wb.Worksheets.Item("CUL")

# Row 36
$ws.Range("H36").Value = 5909.6665
$ws.Range("I36").Value = 1364.5
$ws.Range("K36").Value = 4093.5
$ws.Range("M36").Value = -3924.5

# Row 107
$ws.Range("H107").Value = 405.16666
$ws.Range("I107").Value = 175
$ws.Range("J107").Value = 433.9375
$ws.Range("K107").Value = 525
$ws.Range("L107").Value = 1301.8125
$ws.Range("M107").Value = 1395
$ws.Range("N107").Value = -5141.8125

# Row 128
$ws.Range("H128").Value = 3979888
$ws.Range("I128").Value = 3979888
$ws.Range("K128").Value = 11939664
$ws.Range("M128").Value = -11934684


$ws = $wb.Worksheets.Item("GSM")

# Row 20
$ws.Range("H20").Value = 35317.145
$ws.Range("J20").Value = 37036.668
$ws.Range("L20").Value = 37036.668
$ws.Range("N20").Value = -37526.668

# Row 24
$ws.Range("H24").Value = 2119555.5
$ws.Range("J24").Value = 2243950.5
$ws.Range("L24").Value = 2243950.5
$ws.Range("N24").Value = -2244296.5

# Row 43
$ws.Range("H43").Value = 2553.75
$ws.Range("I43").Value = 1405.6666
$ws.Range("J43").Value = 5998
$ws.Range("K43").Value = 1405.6666
$ws.Range("L43").Value = 5998
$ws.Range("M43").Value = -1254.6666
$ws.Range("N43").Value = -6300

# Row 59
$ws.Range("H59").Value = 18250
$ws.Range("J59").Value = 0
$ws.Range("L59").Value = 0
$ws.Range("N59").ClearContents()

# Row 70
$ws.Range("H70").Value = 7824.8887
$ws.Range("I70").Value = 7481.3335
$ws.Range("K70").Value = 7481.3335
$ws.Range("M70").Value = -7211.3335

# Row 73
$ws.Range("H73").Value = 7824.8887
$ws.Range("I73").Value = 7481.3335
$ws.Range("K73").Value = 7481.3335
$ws.Range("M73").Value = -6545.3335

# Row 122
$ws.Range("H122").Value = 39233
$ws.Range("I122").Value = 2088.9524
$ws.Range("K122").Value = 6266.8572
$ws.Range("M122").Value = -3816.8572

# Row 132
$ws.Range("H132").Value = 3952.3076
$ws.Range("I132").Value = 3436.75
$ws.Range("K132").Value = 10310.25
$ws.Range("M132").Value = -7780.25


$ws = $wb.Worksheets.Item("LTW")

# Row 11
$ws.Range("H11").Value = 2000
$ws.Range("J11").Value = 2000
$ws.Range("L11").Value = 2000
$ws.Range("N11").Value = -2280

# Row 93
$ws.Range("H93").Value = 1948.6
$ws.Range("I93").Value = 1948.6
$ws.Range("K93").Value = 1948.6
$ws.Range("M93").Value = -700.5999999999999

# Row 136
$ws.Range("H136").Value = 3644.5557
$ws.Range("J136").Value = 4098
$ws.Range("L136").Value = 12294
$ws.Range("N136").Value = -17394


$ws = $wb.Worksheets.Item("WVR")

# Row 31
$ws.Range("H31").Value = 20000
$ws.Range("J31").Value = 20000
$ws.Range("L31").Value = 20000
$ws.Range("N31").Value = -20696


